$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for the table (TruckID, AssignedDockPosition, start_loading_time, end_loading_time)
# starting at row 2 through row 13 (rows 12-13 are newly added)
$data = @(
    @(1, 1, 5, 5),
    @(3, 1, 10, 10),
    @(7, 1, 15, 15),
    @(10, 2, 5, 7),
    @(2, 3, 5, 7),
    @(5, 3, 12, 12),
    @(10, 3, 17, 19),
    @(4, 4, 5, 6),
    @(6, 4, 11, 12),
    @(8, 4, 17, 17),
    @(9, 4, 22, 22),
    @(11, 4, 27, 27)
)

$row = 2
foreach ($vals in $data) {
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $row++
}
